# "Refactored text sizing in GUI"
#
# Six pairs of cells on Sheet1 had their contents (shared-text value plus
# the cell's fill/font formatting) swapped with one another. Implement
# each swap by round-tripping through a scratch cell far outside the
# sheet's used range (A1:AF36) with Range.Copy, which carries both the
# value and the full formatting (fill color, font, border, number format)
# of the source range to the destination - exactly like copy/pasting the
# cells in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("AZ100")

function Swap-Cells($ref1, $ref2) {
    $c1 = $ws.Range($ref1)
    $c2 = $ws.Range($ref2)

    $c1.Copy($scratch)
    $c2.Copy($c1)
    $scratch.Copy($c2)
}

Swap-Cells "E3" "E4"
Swap-Cells "X4" "Z4"
Swap-Cells "Y12" "AA12"
Swap-Cells "W21" "W22"
Swap-Cells "D22" "D24"
Swap-Cells "O22" "O24"

# Scratch cell is scratch again - make sure it is left blank.
$scratch.Clear()

# The saved view's active cell/selection moved from Q17 to O24.
$ws.Range("O24").Select()
